$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original Text storage type (they were stored as inline strings,
# many values look numeric and would otherwise be auto-converted to Number by COM).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.15'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.05'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.408'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06002'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.388'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8076'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9287'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1419'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07435'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03386'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03044'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09356'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.938'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '13MCDexMCB'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001600'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04835'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '15CoinExTokenCET'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'One'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0005943'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '16OneONE'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005379'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004158'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009868'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.00007304'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.663'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.423'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03969'
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1075'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002711'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003032'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.006805'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005198'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9805'
